# "last commit of the day"
# Swap the thick-film 0805 resistors for MELF-style resistors (new Mouser
# part numbers) on rows 16/17/24/25, and fill in the previously-blank
# 100nF foil-capacitor line on row 34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: 7k5 resistor -> MELF variant
$ws.Range("B16").Value = "7k5 MELF Widerstand"
$ws.Range("C16").Value = "MELF Widerstände 1/4watt 7.5Kohms 1% 50ppm 13"" "
$ws.Range("F16").Value = "71-SMM02040C7501FB00"

# Row 17: 1k2 resistor -> MELF variant
$ws.Range("B17").Value = "1k2 MELF Widerstand"
$ws.Range("C17").Value = "MELF Widerstände 1/4watt 1.2Kohms 1% 50ppm 13"""
$ws.Range("F17").Value = "71-SMM02040C1201FB00"

# Row 24: 68k resistor -> MELF variant (previously blank details column)
$ws.Range("B24").Value = "68k MELF Widerstand"
$ws.Range("C24").Value = "MELF Widerstände 1/4watt 68Kohms 1% 50ppm 13"""
$ws.Range("F24").Value = "71-SMM02040C6802FB00"

# Row 25: 6k8 resistor -> MELF variant (previously blank details column)
$ws.Range("B25").Value = "6k8 MELF Widerstand"
$ws.Range("C25").Value = "MELF Widerstände 1/4watt 6.8Kohms 1% 50ppm 13"""
$ws.Range("F25").Value = "71-SMM02040C6801FB00"

# Row 34: new line for a 100nF foil capacitor (previously an empty row)
$ws.Range("B34").Value = "100nF Folienkondensator"
$ws.Range("C34").Value = "Folienkondensatoren 0.1uF 16VDC 5% PPS FILM 1210"
$ws.Range("E34").Value = "Mouser"
$ws.Range("F34").Value = "667-ECH-U1C104JX5"

# Row heights: wrapped-text rows whose content changed get re-measured by
# Excel; rows 16/24/25 now wrap to a single line (33pt, matching the other
# one-line rows), row 34 goes from the blank-row default to a one-line 33pt.
$ws.Rows(16).RowHeight = 33
$ws.Rows(24).RowHeight = 33
$ws.Rows(25).RowHeight = 33
$ws.Rows(34).RowHeight = 33

# Restore the view/selection to what was on screen when the file was saved.
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("F31").Select()
